# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Especial / Primera, Frutilla, Provincia de
# Melipilla, 2023-11-28) above the existing data block that starts at row 529,
# pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows before row 529 (shifts 529:638 down to 531:640) ---
$ws.Rows.Item(529).Insert()
$ws.Rows.Item(529).Insert()

# --- Row 529: Especial, $12.000-13.000, Provincia de Melipilla ---
$ws.Range("A529").Value = 7
$ws.Range("B529").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C529").Value = "Ñuble"
$ws.Range("D529").Value = 45258
$ws.Range("E529").Value = 16
$ws.Range("F529").Value = "Fruta"
$ws.Range("G529").Value = 100101
$ws.Range("H529").Value = "Berries"
$ws.Range("I529").Value = 100112025
$ws.Range("J529").Value = "Frutilla"
$ws.Range("K529").Value = "Sin especificar"
$ws.Range("L529").Value = "Especial"
$ws.Range("M529").Value = 120
$ws.Range("N529").Value = 12000
$ws.Range("O529").Value = 13000
$ws.Range("P529").Value = 12500
$ws.Range("Q529").Value = "$/bandeja 7 kilos"
$ws.Range("R529").Value = "Provincia de Melipilla"
$ws.Range("S529").Value = 1786
$ws.Range("T529").Value = 7

# --- Row 530: Primera, $10.000, Provincia de Melipilla ---
$ws.Range("A530").Value = 7
$ws.Range("B530").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C530").Value = "Ñuble"
$ws.Range("D530").Value = 45258
$ws.Range("E530").Value = 16
$ws.Range("F530").Value = "Fruta"
$ws.Range("G530").Value = 100101
$ws.Range("H530").Value = "Berries"
$ws.Range("I530").Value = 100112025
$ws.Range("J530").Value = "Frutilla"
$ws.Range("K530").Value = "Sin especificar"
$ws.Range("L530").Value = "Primera"
$ws.Range("M530").Value = 60
$ws.Range("N530").Value = 10000
$ws.Range("O530").Value = 10000
$ws.Range("P530").Value = 10000
$ws.Range("Q530").Value = "$/bandeja 7 kilos"
$ws.Range("R530").Value = "Provincia de Melipilla"
$ws.Range("S530").Value = 1429
$ws.Range("T530").Value = 7
